$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Non-Accompanied Group" label to "Control Group" wherever it
# appears in the used range (column C, rows 31-51 hold this group label).
$used = $ws.UsedRange
$rowCount = $used.Rows.Count()
$colCount = $used.Columns.Count()

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -eq "Non-Accompanied Group") {
            $cell.Value = "Control Group"
        }
    }
}

# Update the saved view state (zoom level & active selection) to match the
# state captured when the workbook was last edited.
$excel.ActiveWindow.Zoom = 177
[void]$ws.Range("M40").Select()
